$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for rows 162-201 (A, B, C, D, E, F, G)
$data = @(
    @(2598, 24999, 1, "[3158]", 3158, 3941, 0.8013),
    @(2599, 224999, 1, "[3257]", 3257, 4118, 0.7909),
    @(2600, 1649999, 1, "[3223]", 3223, 4024, 0.8008999999999999),
    @(2601, 199999, 1, "[3208]", 3208, 4043, 0.7935),
    @(2602, 924999, 1, "[3233]", 3233, 4027, 0.8028),
    @(2603, 99999, 1, "[3280]", 3280, 4103, 0.7994),
    @(2604, 24999, 1, "[3237]", 3237, 4037, 0.8018),
    @(2605, 349999, 1, "[3142]", 3142, 3947, 0.796),
    @(2606, 99999, 1, "[3129]", 3129, 3946, 0.793),
    @(2607, 24999, 1, "[3126]", 3126, 3891, 0.8034),
    @(2608, 49999, 1, "[3203]", 3203, 4030, 0.7948),
    @(2609, 1024999, 1, "[3240]", 3240, 4038, 0.8024),
    @(2610, 149999, 1, "[3302]", 3302, 4101, 0.8052),
    @(2611, 324999, 1, "[3269]", 3269, 4086, 0.8),
    @(2612, 49999, 1, "[3191]", 3191, 3986, 0.8006),
    @(2613, 99999, 1, "[3154]", 3154, 3936, 0.8013),
    @(2614, 49999, 1, "[3176]", 3176, 3971, 0.7998),
    @(2615, 49999, 1, "[3209]", 3209, 4007, 0.8008),
    @(2616, 199999, 1, "[3316]", 3316, 4080, 0.8127),
    @(2617, 49999, 1, "[3079]", 3079, 3918, 0.7859),
    @(2618, 274999, 1, "[3213]", 3213, 4012, 0.8008),
    @(2619, 99999, 1, "[3239]", 3239, 4015, 0.8067),
    @(2620, 99999, 1, "[3216]", 3216, 3971, 0.8099),
    @(2621, 74999, 1, "[3208]", 3208, 3990, 0.804),
    @(2622, 74999, 1, "[3170]", 3170, 3951, 0.8023),
    @(2623, 524999, 1, "[3234]", 3234, 4020, 0.8045),
    @(2624, 74999, 1, "[3251]", 3251, 4053, 0.8021),
    @(2625, 24999, 1, "[3140]", 3140, 3940, 0.797),
    @(2626, 49999, 1, "[3193]", 3193, 3984, 0.8015),
    @(2627, 74999, 1, "[3152]", 3152, 3961, 0.7958),
    @(2628, 24999, 1, "[3255]", 3255, 3996, 0.8146),
    @(2629, 774999, 1, "[3175]", 3175, 3909, 0.8122),
    @(2630, 74999, 1, "[3217]", 3217, 4038, 0.7967),
    @(2631, 449999, 1, "[3239]", 3239, 4055, 0.7988),
    @(2632, 49999, 1, "[3131]", 3131, 3981, 0.7865),
    @(2633, 99999, 1, "[3263]", 3263, 4066, 0.8025),
    @(2634, 924999, 1, "[3150]", 3150, 3960, 0.7955),
    @(2635, 24999, 1, "[3099]", 3099, 3924, 0.7897999999999999),
    @(2636, 24999, 1, "[3211]", 3211, 4055, 0.7919),
    @(2637, 149999, 1, "[3228]", 3228, 3947, 0.8178)
)

$startRow = 162
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}
